$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
# B4 held the text "Next" (style matching B/C "next" marker). That marker
# moves to the new D4 cell, and B4 becomes a plain date using the same
# style C4/B3 already use.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Value = "Next"

$ws.Range("C4").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = 45910

# S4 gains a date value (style unchanged)
$ws.Range("S4").Value = 45912

# --- Row 5 ---
# T5 keeps its value but its style is normalized to match the rest of
# column T (same style as T4/T6/T7)
$ws.Range("T4").Copy() | Out-Null
$ws.Range("T5").PasteSpecial(-4122) | Out-Null
$ws.Range("T5").Value = 45903

# --- Row 6 ---
# J6 becomes a "Next" marker cell, matching the style already used for
# that marker in the J/K/L columns (as seen on L7 before this edit).
$ws.Range("L7").Copy() | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").Value = "Next"

# --- Row 7 ---
# L7 stops being the "Next" marker and becomes a normal date, matching
# the style used elsewhere in column L (e.g. L3/L4/L5).
$ws.Range("L5").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null
$ws.Range("L7").Value = 45917

# --- Row 8 ---
# A8 gains a date value (style unchanged)
$ws.Range("A8").Value = 45909
# P8 gains a date value, matching the style used in column P (e.g. P6)
$ws.Range("P6").Copy() | Out-Null
$ws.Range("P8").PasteSpecial(-4122) | Out-Null
$ws.Range("P8").Value = 45910

# --- Row 11 (new row) ---
# I11 gains a date value, matching the style used elsewhere in column I
$ws.Range("I10").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").Value = 45917

$excel.CutCopyMode = 0

# --- Sheet view adjustments ---
# Move the selection to match the saved view state in the edited workbook
# (this also naturally drops the stale topLeftCell scroll position).
$ws.Range("U9").Select() | Out-Null
